# Updates crypto "cryptos.xlsx" sheet to match the refreshed symbol-list
# snapshot (GitHub Actions run on Fri Dec 23 22:51:36 UTC 2022).
#
# Column D holds prices as plain text (e.g. "245.98") so that values
# such as "6.380" or "4.470" keep their trailing zero instead of being
# normalised by Excels General number handling. To write a new price
# without Excel silently re-typing the cell as a Number, we temporarily
# flip the cell to a Text format before assigning the value, then copy
# the plain (un-styled) look of a neighboring, unchanged price cell
# (D3) back onto it so no stray formatting sticks around.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("D3").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.02"
$ws.Range("D2").Style = $plainStyle
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.299"
$ws.Range("D4").Style = $plainStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05874"
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.385"
$ws.Range("D6").Style = $plainStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.380"
$ws.Range("D7").Style = $plainStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8165"
$ws.Range("D8").Style = $plainStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9612"
$ws.Range("D9").Style = $plainStyle
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1418"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03595"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07322"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03048"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.470"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09389"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001606"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04811"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005905"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006246"
$ws.Range("D19").Style = $plainStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009875"
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009709"
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.683"
$ws.Range("D23").Style = $plainStyle
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3266"
$ws.Range("D25").Style = $plainStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1276"
$ws.Range("D26").Style = $plainStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002474"
$ws.Range("D27").Style = $plainStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03873"
$ws.Range("D40").Style = $plainStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006626"
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D42").Style = $plainStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003003"
$ws.Range("D43").Style = $plainStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005902"
$ws.Range("D44").Style = $plainStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D46").Style = $plainStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7757"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04008"
$ws.Range("D48").Style = $plainStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D49").Style = $plainStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
$ws.Range("D50").Style = $plainStyle
